$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A93").Value = "2024-09-29 03:28:41"
$ws.Range("B93").Value = "monitor_price"
$ws.Range("C93").Value = "https://example.com/product"
$ws.Range("D93").Value = "$199.99"
$ws.Range("E93").Value = "2024-09-29"
$ws.Range("F93").Value = "03:28:41"

$ws.Range("A94").Value = "2024-09-29 03:28:41"
$ws.Range("B94").Value = "monitor_price"
$ws.Range("C94").Value = "invalid_url"
$ws.Range("D94").Value = "Error fetching price: Invalid URL"
$ws.Range("E94").Value = "2024-09-29"
$ws.Range("F94").Value = "03:28:41"

$ws.Range("A95").Value = "2024-09-29 03:28:43"
$ws.Range("B95").Value = "monitor_price"
$ws.Range("C95").Value = "https://example.com/product"
$ws.Range("D95").Value = "100 USD"
$ws.Range("E95").Value = "2024-09-29"
$ws.Range("F95").Value = "03:28:43"

$ws.Range("A96").Value = "2024-09-29 03:32:54"
$ws.Range("B96").Value = "monitor_price"
$ws.Range("C96").Value = "https://example.com/product"
$ws.Range("D96").Value = "$199.99"
$ws.Range("E96").Value = "2024-09-29"
$ws.Range("F96").Value = "03:32:54"

$ws.Range("A97").Value = "2024-09-29 03:32:54"
$ws.Range("B97").Value = "monitor_price"
$ws.Range("C97").Value = "invalid_url"
$ws.Range("D97").Value = "Error fetching price: Invalid URL"
$ws.Range("E97").Value = "2024-09-29"
$ws.Range("F97").Value = "03:32:54"

$ws.Range("A98").Value = "2024-09-29 03:32:56"
$ws.Range("B98").Value = "monitor_price"
$ws.Range("C98").Value = "https://example.com/product"
$ws.Range("D98").Value = "100 USD"
$ws.Range("E98").Value = "2024-09-29"
$ws.Range("F98").Value = "03:32:56"

$ws.Range("A99").Value = "2024-09-29 03:48:56"
$ws.Range("B99").Value = "monitor_price"
$ws.Range("C99").Value = "https://example.com/product"
$ws.Range("D99").Value = "$199.99"
$ws.Range("E99").Value = "2024-09-29"
$ws.Range("F99").Value = "03:48:56"

$ws.Range("A100").Value = "2024-09-29 03:48:56"
$ws.Range("B100").Value = "monitor_price"
$ws.Range("C100").Value = "invalid_url"
$ws.Range("D100").Value = "Error fetching price: Invalid URL"
$ws.Range("E100").Value = "2024-09-29"
$ws.Range("F100").Value = "03:48:56"

$ws.Range("A101").Value = "2024-09-29 03:48:58"
$ws.Range("B101").Value = "monitor_price"
$ws.Range("C101").Value = "https://example.com/product"
$ws.Range("D101").Value = "100 USD"
$ws.Range("E101").Value = "2024-09-29"
$ws.Range("F101").Value = "03:48:58"

$ws.Range("A102").Value = "2024-09-29 03:49:43"
$ws.Range("B102").Value = "monitor_price"
$ws.Range("C102").Value = "https://example.com/product"
$ws.Range("D102").Value = "$199.99"
$ws.Range("E102").Value = "2024-09-29"
$ws.Range("F102").Value = "03:49:43"

$ws.Range("A103").Value = "2024-09-29 03:49:43"
$ws.Range("B103").Value = "monitor_price"
$ws.Range("C103").Value = "invalid_url"
$ws.Range("D103").Value = "Error fetching price: Invalid URL"
$ws.Range("E103").Value = "2024-09-29"
$ws.Range("F103").Value = "03:49:43"

$ws.Range("A104").Value = "2024-09-29 03:49:45"
$ws.Range("B104").Value = "monitor_price"
$ws.Range("C104").Value = "https://example.com/product"
$ws.Range("D104").Value = "100 USD"
$ws.Range("E104").Value = "2024-09-29"
$ws.Range("F104").Value = "03:49:45"

$ws.Range("A105").Value = "2024-09-29 03:52:26"
$ws.Range("B105").Value = "monitor_price"
$ws.Range("C105").Value = "https://example.com/product"
$ws.Range("D105").Value = "$199.99"
$ws.Range("E105").Value = "2024-09-29"
$ws.Range("F105").Value = "03:52:26"

$ws.Range("A106").Value = "2024-09-29 03:52:27"
$ws.Range("B106").Value = "monitor_price"
$ws.Range("C106").Value = "invalid_url"
$ws.Range("D106").Value = "Error fetching price: Invalid URL"
$ws.Range("E106").Value = "2024-09-29"
$ws.Range("F106").Value = "03:52:27"

$ws.Range("A107").Value = "2024-09-29 03:52:29"
$ws.Range("B107").Value = "monitor_price"
$ws.Range("C107").Value = "https://example.com/product"
$ws.Range("D107").Value = "100 USD"
$ws.Range("E107").Value = "2024-09-29"
$ws.Range("F107").Value = "03:52:29"

$ws.Range("A108").Value = "2024-09-29 03:53:03"
$ws.Range("B108").Value = "monitor_price"
$ws.Range("C108").Value = "https://example.com/product"
$ws.Range("D108").Value = "$199.99"
$ws.Range("E108").Value = "2024-09-29"
$ws.Range("F108").Value = "03:53:03"

$ws.Range("A109").Value = "2024-09-29 03:53:03"
$ws.Range("B109").Value = "monitor_price"
$ws.Range("C109").Value = "invalid_url"
$ws.Range("D109").Value = "Error fetching price: Invalid URL"
$ws.Range("E109").Value = "2024-09-29"
$ws.Range("F109").Value = "03:53:03"

$ws.Range("A110").Value = "2024-09-29 03:53:05"
$ws.Range("B110").Value = "monitor_price"
$ws.Range("C110").Value = "https://example.com/product"
$ws.Range("D110").Value = "100 USD"
$ws.Range("E110").Value = "2024-09-29"
$ws.Range("F110").Value = "03:53:05"

